$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '332.64'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.99%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '45.51'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.62%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.528'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.37%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08490'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '5.86%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.076'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.21%'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9884'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.80%'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.541'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-5.47%'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1159'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.33%'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1919'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.22%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '9.506'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-6.52%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09783'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.14%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04690'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-3.94%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1061'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.20%'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001276'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.34%'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005903'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.23%'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.387'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.17%'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.434'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.56%'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3355'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.44%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1384'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.02%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2553'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.10%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04162'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.73%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001302'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.05%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004600'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '5.46%'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '10.42%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0002987'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-20.36%'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02718'

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '5.86%'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05749'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-0.14%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007769'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.40%'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1434'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.32%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007630'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.04%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002088'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.66%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008054'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-11.06%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3557'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007063'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.65%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000752'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.10%'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '0.13%'

$ws.Range("B49").Value = 'CoinbaseStockToken'

$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003540'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '1.02%'

$ws.Range("B50").Value = 'BOLO'

$ws.Range("C50").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.003386'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-3.34%'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00002106'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.10%'
